$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "总计" (summary) sheet: insert a new row for "2022-Q4" right under
#    the header, pushing the existing quarters down by one row and
#    renumbering the index column (A) sequentially.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows("2:2").Insert()
$summary.Range("B2:D2").ClearFormats()

# Give the new A2 the same formatting as the other index cells in column A.
$summary.Range("A3").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("B2").Value = "2022-Q4"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0

$summary.Range("A3").Value = 1
$summary.Range("B3").Value = "2022-Q1"
$summary.Range("C3").Value = 5
$summary.Range("D3").Value = 0.07000000000000001

$summary.Range("A4").Value = 2
$summary.Range("B4").Value = "2021-Q4"
$summary.Range("C4").Value = 2
$summary.Range("D4").Value = 0.01

$summary.Range("A5").Value = 3
$summary.Range("B5").Value = "2021-Q3"
$summary.Range("C5").Value = 3
$summary.Range("D5").Value = 0.11

$summary.Range("A6").Value = 4
$summary.Range("B6").Value = "2020-Q4"
$summary.Range("C6").Value = 2
$summary.Range("D6").Value = 0.02

# ---------------------------------------------------------------------
# 2) Add the new "2022-Q4" detail sheet as the 2nd sheet (right after
#    "总计"), with the same layout as the other quarterly sheets.
# ---------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $firstSheet)
$newSheet.Name = "2022-Q4"

# Match the look & feel (page setup / outline options) used by the other
# data sheets in this workbook.
$newSheet.Outline.SummaryRow = 1
$newSheet.Outline.SummaryColumn = 1
$newSheet.PageSetup.LeftMargin = 54
$newSheet.PageSetup.RightMargin = 54
$newSheet.PageSetup.TopMargin = 72
$newSheet.PageSetup.BottomMargin = 72
$newSheet.PageSetup.HeaderMargin = 36
$newSheet.PageSetup.FooterMargin = 36

# Header row.
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Copy the header formatting (bold font + border + centered alignment,
# style index "2" in the original sheets) from the summary sheet.
$summary.Range("B1:D1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# Data row. Numeric-looking text (fund code / percentages) is entered with
# a leading apostrophe so Excel keeps it as text (preserving leading /
# trailing zeros) instead of converting it to a number.
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "'005120"
$newSheet.Range("C2").Value = "上投摩根量化多因子灵活配置混合"
$newSheet.Range("D2").Value = "'0.19"
$newSheet.Range("E2").Value = "'94.61"
$newSheet.Range("F2").Value = "'1.59"
$newSheet.Range("G2").Value = "'0.0030"
$newSheet.Range("H2").Value = 6

# Give A2 the same style as A1 (bold/border/center), matching the "index"
# column style used throughout the workbook.
$summary.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$newSheet.Range("A2").Value = 0

# ---------------------------------------------------------------------
# 3) Re-activate the last sheet ("2020-Q4") so it stays the selected tab,
#    as it was before this edit (adding a sheet makes it active by
#    default).
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$lastSheet.Activate()
